$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order ids refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961321501656"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961338068857"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961338068857"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961338548841"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961339268878"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961321181524.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961321341662.csv"
$ws1.Range("B4").Value = "go_stims-16509961321341662.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961321501656.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509961337908444.csv"
$ws2.Range("B3").Value = "ZB-match_4-16509961323501673.csv"
$ws2.Range("B4").Value = "TB-16509961333428776.csv"
$ws2.Range("B5").Value = "ZB-match_8-165099613250288.csv"
$ws2.Range("B6").Value = "OB-16509961327428467.csv"
$ws2.Range("B7").Value = "OB-16509961332228456.csv"
$ws2.Range("B8").Value = "ZB-match_6-16509961325428438.csv"
$ws2.Range("B9").Value = "TB-16509961334388468.csv"
$ws2.Range("B10").Value = "OB-16509961327108471.csv"

# --- Sheet 3 (RS) --- (no cell content changes, only the name change above)

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961338228457.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961338068857.csv"
$ws4.Range("B4").Value = "MM_stims-16509961338388793.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961338228457.csv"
$ws4.Range("B6").Value = "MM_stims-16509961338548841.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961338388793.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961338708506.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961338948429.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650996133910885.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961338548841.csv"
